$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.273.04'
$ws.Range('E2').Value = '  +3.01%  '

$ws.Range('D3').Value = '1.813.63'
$ws.Range('E3').Value = '  +0.94%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.25%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '339.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.25%  '

$ws.Range('E6').Value = '  -0.07%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3910'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.77%  '

$ws.Range('E8').Value = '  +0.69%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.40'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.64%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.192'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.95%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07569'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.55%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9987'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.26%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.12'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.16%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.511'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.50%  '

$ws.Range('D15').Value = '1.813.54'
$ws.Range('E15').Value = '  +0.95%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.142'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.79%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001103'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.37%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06691'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.32%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '84.96'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9994'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.78'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.00%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.560'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.24%  '

$ws.Range('D23').Value = '28.227.51'
$ws.Range('E23').Value = '  +2.88%  '

$ws.Range('E24').Value = '  -0.74%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.401'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.82%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.490'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.24%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.528'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.59%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.88%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '153.76'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.57%  '

$ws.Range('D30').Value = '2.019.55'
$ws.Range('E30').Value = '  +0.94%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '135.35'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.151'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.55%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.019'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.70%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08824'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.26%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.01%  '

$ws.Range('B36').Value = 'TheSandbox'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6960'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.46%  '

$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.469'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.02%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02424'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.21%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06538'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.13%  '

$ws.Range('E40').Value = '  -2.56%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2213'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.26%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.259'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.16%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.483'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.65%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.20%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6436'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.51%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.870'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.13%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.153'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.65%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.73'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.94%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07198'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.05%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '80.03'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.28%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.246'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.10%  '
